$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (outside the used A:G range) used to stage each new value as
# genuine text (NumberFormat "@") and then Copy/PasteSpecial values-only into
# the destination cell. This preserves literal text like "309.17", "0.1160",
# "1.05%" or "0" exactly as text instead of letting Excel's normal
# autoconvert-on-entry turn them into numbers/percentages, and it avoids
# minting a new cell style (NumberFormat assigned directly to the destination
# cell would fork its style away from the original).
$helper = $ws.Cells.Item(1, 10)
$helper.NumberFormat = "@"

$helper.Value = '309.17'
$helper.Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4163)
$helper.Value = '0.74%'
$helper.Copy()
$ws.Cells.Item(2, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(2, 7).PasteSpecial(-4163)
$helper.Value = '41.21'
$helper.Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4163)
$helper.Value = '0.04%'
$helper.Copy()
$ws.Cells.Item(3, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(3, 7).PasteSpecial(-4163)
$helper.Value = '5.253'
$helper.Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4163)
$helper.Value = '2.34%'
$helper.Copy()
$ws.Cells.Item(4, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(4, 7).PasteSpecial(-4163)
$helper.Value = '0.07664'
$helper.Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4163)
$helper.Value = '0.95%'
$helper.Copy()
$ws.Cells.Item(5, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(5, 7).PasteSpecial(-4163)
$helper.Value = 'GateToken'
$helper.Copy()
$ws.Cells.Item(6, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$helper.Copy()
$ws.Cells.Item(6, 3).PasteSpecial(-4163)
$helper.Value = '4.336'
$helper.Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)
$helper.Value = '1.66%'
$helper.Copy()
$ws.Cells.Item(6, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(6, 7).PasteSpecial(-4163)
$helper.Value = 'FTXToken'
$helper.Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$helper.Copy()
$ws.Cells.Item(7, 3).PasteSpecial(-4163)
$helper.Value = '1.618'
$helper.Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4163)
$helper.Value = '0.72%'
$helper.Copy()
$ws.Cells.Item(7, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(7, 7).PasteSpecial(-4163)
$helper.Value = 'MXToken'
$helper.Copy()
$ws.Cells.Item(8, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$helper.Copy()
$ws.Cells.Item(8, 3).PasteSpecial(-4163)
$helper.Value = '0.9182'
$helper.Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4163)
$helper.Value = '1.90%'
$helper.Copy()
$ws.Cells.Item(8, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(8, 7).PasteSpecial(-4163)
$helper.Value = 'BTSEToken'
$helper.Copy()
$ws.Cells.Item(9, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$helper.Copy()
$ws.Cells.Item(9, 3).PasteSpecial(-4163)
$helper.Value = '2.445'
$helper.Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4163)
$helper.Value = '0.58%'
$helper.Copy()
$ws.Cells.Item(9, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(9, 7).PasteSpecial(-4163)
$helper.Value = 'LiechtensteinCryptoassetsExchange'
$helper.Copy()
$ws.Cells.Item(10, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$helper.Copy()
$ws.Cells.Item(10, 3).PasteSpecial(-4163)
$helper.Value = '0.1222'
$helper.Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4163)
$helper.Value = '12.12%'
$helper.Copy()
$ws.Cells.Item(10, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(10, 7).PasteSpecial(-4163)
$helper.Value = 'WazirX'
$helper.Copy()
$ws.Cells.Item(11, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$helper.Copy()
$ws.Cells.Item(11, 3).PasteSpecial(-4163)
$helper.Value = '0.1833'
$helper.Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4163)
$helper.Value = '4.83%'
$helper.Copy()
$ws.Cells.Item(11, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(11, 7).PasteSpecial(-4163)
$helper.Value = 'MandalaExchangeToken'
$helper.Copy()
$ws.Cells.Item(12, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$helper.Copy()
$ws.Cells.Item(12, 3).PasteSpecial(-4163)
$helper.Value = '0.09133'
$helper.Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)
$helper.Value = '0.08%'
$helper.Copy()
$ws.Cells.Item(12, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(12, 7).PasteSpecial(-4163)
$helper.Value = 'BitrueCoin'
$helper.Copy()
$ws.Cells.Item(13, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$helper.Copy()
$ws.Cells.Item(13, 3).PasteSpecial(-4163)
$helper.Value = '0.04317'
$helper.Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4163)
$helper.Value = '-0.03%'
$helper.Copy()
$ws.Cells.Item(13, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(13, 7).PasteSpecial(-4163)
$helper.Value = 'BitMartToken'
$helper.Copy()
$ws.Cells.Item(14, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$helper.Copy()
$ws.Cells.Item(14, 3).PasteSpecial(-4163)
$helper.Value = '0.1052'
$helper.Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4163)
$helper.Value = '0.08%'
$helper.Copy()
$ws.Cells.Item(14, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(14, 7).PasteSpecial(-4163)
$helper.Value = 'BitForexToken'
$helper.Copy()
$ws.Cells.Item(15, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$helper.Copy()
$ws.Cells.Item(15, 3).PasteSpecial(-4163)
$helper.Value = '0.001261'
$helper.Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4163)
$helper.Value = '0.33%'
$helper.Copy()
$ws.Cells.Item(15, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(15, 7).PasteSpecial(-4163)
$helper.Value = 'TigerCash'
$helper.Copy()
$ws.Cells.Item(16, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$helper.Copy()
$ws.Cells.Item(16, 3).PasteSpecial(-4163)
$helper.Value = '0.005795'
$helper.Copy()
$ws.Cells.Item(16, 4).PasteSpecial(-4163)
$helper.Value = '-1.96%'
$helper.Copy()
$ws.Cells.Item(16, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(16, 7).PasteSpecial(-4163)
$helper.Value = 'UpBots'
$helper.Copy()
$ws.Cells.Item(17, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$helper.Copy()
$ws.Cells.Item(17, 3).PasteSpecial(-4163)
$helper.Value = '0.007498'
$helper.Copy()
$ws.Cells.Item(17, 4).PasteSpecial(-4163)
$helper.Value = '2,389.02%'
$helper.Copy()
$ws.Cells.Item(17, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(17, 7).PasteSpecial(-4163)
$helper.Value = 'LEO'
$helper.Copy()
$ws.Cells.Item(18, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$helper.Copy()
$ws.Cells.Item(18, 3).PasteSpecial(-4163)
$helper.Value = '3.347'
$helper.Copy()
$ws.Cells.Item(18, 4).PasteSpecial(-4163)
$helper.Value = '-0.20%'
$helper.Copy()
$ws.Cells.Item(18, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(18, 7).PasteSpecial(-4163)
$helper.Value = 'BitpandaEcosystemToken'
$helper.Copy()
$ws.Cells.Item(19, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$helper.Copy()
$ws.Cells.Item(19, 3).PasteSpecial(-4163)
$helper.Value = '0.3335'
$helper.Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4163)
$helper.Value = '1.85%'
$helper.Copy()
$ws.Cells.Item(19, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(19, 7).PasteSpecial(-4163)
$helper.Value = 'MCDex'
$helper.Copy()
$ws.Cells.Item(20, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$helper.Copy()
$ws.Cells.Item(20, 3).PasteSpecial(-4163)
$helper.Value = '7.281'
$helper.Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4163)
$helper.Value = '11.27%'
$helper.Copy()
$ws.Cells.Item(20, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(20, 7).PasteSpecial(-4163)
$helper.Value = 'ProBitToken'
$helper.Copy()
$ws.Cells.Item(21, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$helper.Copy()
$ws.Cells.Item(21, 3).PasteSpecial(-4163)
$helper.Value = '0.1402'
$helper.Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4163)
$helper.Value = '2.68%'
$helper.Copy()
$ws.Cells.Item(21, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(21, 7).PasteSpecial(-4163)
$helper.Value = 'ZBToken'
$helper.Copy()
$ws.Cells.Item(22, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$helper.Copy()
$ws.Cells.Item(22, 3).PasteSpecial(-4163)
$helper.Value = '0.2916'
$helper.Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4163)
$helper.Value = '8.61%'
$helper.Copy()
$ws.Cells.Item(22, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(22, 7).PasteSpecial(-4163)
$helper.Value = 'CoinExToken'
$helper.Copy()
$ws.Cells.Item(23, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$helper.Copy()
$ws.Cells.Item(23, 3).PasteSpecial(-4163)
$helper.Value = '0.04066'
$helper.Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4163)
$helper.Value = '-2.46%'
$helper.Copy()
$ws.Cells.Item(23, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(23, 7).PasteSpecial(-4163)
$helper.Value = 'BitKan'
$helper.Copy()
$ws.Cells.Item(24, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$helper.Copy()
$ws.Cells.Item(24, 3).PasteSpecial(-4163)
$helper.Value = '0.001261'
$helper.Copy()
$ws.Cells.Item(24, 4).PasteSpecial(-4163)
$helper.Value = '3.01%'
$helper.Copy()
$ws.Cells.Item(24, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(24, 7).PasteSpecial(-4163)
$helper.Value = 'HotbitToken'
$helper.Copy()
$ws.Cells.Item(25, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$helper.Copy()
$ws.Cells.Item(25, 3).PasteSpecial(-4163)
$helper.Value = '0.004083'
$helper.Copy()
$ws.Cells.Item(25, 4).PasteSpecial(-4163)
$helper.Value = '-0.11%'
$helper.Copy()
$ws.Cells.Item(25, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(25, 7).PasteSpecial(-4163)
$helper.Value = 'NitroEx'
$helper.Copy()
$ws.Cells.Item(26, 2).PasteSpecial(-4163)
$helper.Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$helper.Copy()
$ws.Cells.Item(26, 3).PasteSpecial(-4163)
$helper.Value = '0.0001273'
$helper.Copy()
$ws.Cells.Item(26, 4).PasteSpecial(-4163)
$helper.Value = '-2.19%'
$helper.Copy()
$ws.Cells.Item(26, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(26, 7).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(27, 7).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(28, 7).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(29, 7).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(30, 7).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(31, 7).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(32, 7).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(33, 7).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(34, 7).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(35, 7).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(36, 7).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(37, 7).PasteSpecial(-4163)
$helper.Value = '0.02440'
$helper.Copy()
$ws.Cells.Item(38, 4).PasteSpecial(-4163)
$helper.Value = '4.17%'
$helper.Copy()
$ws.Cells.Item(38, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(38, 7).PasteSpecial(-4163)
$helper.Value = '0.05274'
$helper.Copy()
$ws.Cells.Item(39, 4).PasteSpecial(-4163)
$helper.Value = '2.52%'
$helper.Copy()
$ws.Cells.Item(39, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(39, 7).PasteSpecial(-4163)
$helper.Value = '0.007855'
$helper.Copy()
$ws.Cells.Item(40, 4).PasteSpecial(-4163)
$helper.Value = '1.09%'
$helper.Copy()
$ws.Cells.Item(40, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(40, 7).PasteSpecial(-4163)
$helper.Value = '0.1313'
$helper.Copy()
$ws.Cells.Item(41, 4).PasteSpecial(-4163)
$helper.Value = '1.30%'
$helper.Copy()
$ws.Cells.Item(41, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(41, 7).PasteSpecial(-4163)
$helper.Value = '0.006825'
$helper.Copy()
$ws.Cells.Item(42, 4).PasteSpecial(-4163)
$helper.Value = '-2.03%'
$helper.Copy()
$ws.Cells.Item(42, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(42, 7).PasteSpecial(-4163)
$helper.Value = '0.001914'
$helper.Copy()
$ws.Cells.Item(43, 4).PasteSpecial(-4163)
$helper.Value = '-2.92%'
$helper.Copy()
$ws.Cells.Item(43, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(43, 7).PasteSpecial(-4163)
$helper.Value = '0.008348'
$helper.Copy()
$ws.Cells.Item(44, 4).PasteSpecial(-4163)
$helper.Value = '-2.18%'
$helper.Copy()
$ws.Cells.Item(44, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(44, 7).PasteSpecial(-4163)
$helper.Value = '0.3349'
$helper.Copy()
$ws.Cells.Item(45, 4).PasteSpecial(-4163)
$helper.Value = '10.02%'
$helper.Copy()
$ws.Cells.Item(45, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(45, 7).PasteSpecial(-4163)
$helper.Value = '0.00006846'
$helper.Copy()
$ws.Cells.Item(46, 4).PasteSpecial(-4163)
$helper.Value = '6.81%'
$helper.Copy()
$ws.Cells.Item(46, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(46, 7).PasteSpecial(-4163)
$helper.Value = '0.00000000752'
$helper.Copy()
$ws.Cells.Item(47, 4).PasteSpecial(-4163)
$helper.Value = '0.02%'
$helper.Copy()
$ws.Cells.Item(47, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(47, 7).PasteSpecial(-4163)
$helper.Value = '1,971.62%'
$helper.Copy()
$ws.Cells.Item(48, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(48, 7).PasteSpecial(-4163)
$helper.Value = '-6.87%'
$helper.Copy()
$ws.Cells.Item(49, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(49, 7).PasteSpecial(-4163)
$helper.Value = '0.00002104'
$helper.Copy()
$ws.Cells.Item(50, 4).PasteSpecial(-4163)
$helper.Value = '0.02%'
$helper.Copy()
$ws.Cells.Item(50, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(50, 7).PasteSpecial(-4163)
$helper.Value = '0.0002004'
$helper.Copy()
$ws.Cells.Item(51, 4).PasteSpecial(-4163)
$helper.Value = '0.02%'
$helper.Copy()
$ws.Cells.Item(51, 5).PasteSpecial(-4163)
$helper.Value = '2'
$helper.Copy()
$ws.Cells.Item(51, 7).PasteSpecial(-4163)

$helper.Clear()
$excel.CutCopyMode = $false
Write-Host "Applied 12-2-2023 symbol list update"
